# Update the PSSM value matrix (B2:K21) with the supplemental-figure
# recomputed scores. Row/column headers (row 1, column A) are unchanged.
$ws = $excel.ActiveWorkbook.ActiveSheet

$newScores = @(
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,2.472395717281631,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,1.776552230007785,-17.81574602556542,2.617216846838455,-17.81574602556542,1.720313326974988,-17.81574602556542,2.502369753609742,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,2.031429104973347,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(2.950987487992853,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,2.082402145832663,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(3.616727060210178,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,1.549126019486971,-17.81574602556542,2.219113373957953),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,1.950389765993019,-17.81574602556542,2.461220054682864,-17.81574602556542,-17.81574602556542,-17.81574602556542,1.286365527500274),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,1.556883139111396,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,2.291150407551088,1.537411571160619),
        @(-17.81574602556542,-17.81574602556542,1.014227499659732,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,2.030388597962802),
        @(-17.81574602556542,-17.81574602556542,0.8828535222414425,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,2.269317910258049,-17.81574602556542),
        @(-17.81574602556542,4.321922154358783,-2.143403561378859,-17.81574602556542,-17.81574602556542,-17.81574602556542,0.4513561067434353,0.6736041355153138,1.10696902488401,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,-17.81574602556542,0.2239709308357506,0.6590254041306797,1.301260338840465,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,2.471381267604031,-17.81574602556542,-17.81574602556542,-17.81574602556542,1.807001127631635,2.123249144962508,-17.81574602556542,-17.81574602556542),
        @(-17.81574602556542,-17.81574602556542,2.794508365308346,-17.81574602556542,3.793299665419652,-17.81574602556542,2.24076071711868,1.997504389185702,-17.81574602556542,2.564830239818593),
        @(-17.81574602556542,-17.81574602556542,-17.81574602556542,3.16315222569113,-17.81574602556542,3.379206159749757,2.571430942996823,-17.81574602556542,-17.81574602556542,-17.81574602556542)
    )

$startRow = 2
$startCol = 2  # column B

for ($r = 0; $r -lt $newScores.Length; $r++) {
    $rowValues = $newScores[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowValues[$c]
    }
}
